$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.917.38'
$ws.Range('E2').Value = '  +1.31%  '
$ws.Range('D3').Value = '1.638.31'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.58'
$ws.Range('E5').Value = '  +0.64%  '
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.62'
$ws.Range('E8').Value = '  +1.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.261'
$ws.Range('E9').Value = '  -0.97%  '
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0875'
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('D12').Value = '1.871.47'
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('D13').Value = '1.647.10'
$ws.Range('E13').Value = '  +1.21%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.576'
$ws.Range('E14').Value = '  +4.06%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.10'
$ws.Range('E15').Value = '  +1.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.97'
$ws.Range('E16').Value = '  +0.87%  '
$ws.Range('D17').Value = '27.910.14'
$ws.Range('E17').Value = '  +1.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '231.94'
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('D19').Value = '0.0₃0724'
$ws.Range('E19').Value = '  +0.79%  '
$ws.Range('E20').Value = '  +0.47%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.79'
$ws.Range('E22').Value = '  +2.41%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  -3.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.69'
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.74'
$ws.Range('E27').Value = '  +1.33%  '
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('E32').Value = '  +1.79%  '
$ws.Range('E33').Value = '  +1.43%  '
$ws.Range('D34').Value = '1.410.28'
$ws.Range('E34').Value = '  -3.89%  '
$ws.Range('E35').Value = '  +1.20%  '
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.891'
$ws.Range('E37').Value = '  +1.33%  '
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.915'
$ws.Range('E40').Value = '  -2.82%  '
$ws.Range('E41').Value = '  -1.10%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.84'
$ws.Range('E43').Value = '  +4.50%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '66.29'
$ws.Range('E44').Value = '  -2.30%  '
$ws.Range('E45').Value = '  +1.46%  '
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('D47').Value = '1.780.23'
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.32'
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('E50').Value = '  +0.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.63'
$ws.Range('E51').Value = '  -0.48%  '
